$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.188.11'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').Value = '1.824.71'
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = "'236.03"
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('D6').Value = "'0.6100"
$ws.Range('E6').Value = '  -2.95%  '
$ws.Range('D7').Value = "'1.002"
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('E8').Value = '  -4.55%  '
$ws.Range('D9').Value = "'0.2804"
$ws.Range('E9').Value = '  -3.00%  '
$ws.Range('D10').Value = "'23.49"
$ws.Range('D11').Value = "'0.07661"
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('D12').Value = '1.826.10'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').Value = "'4.804"
$ws.Range('E13').Value = '  -3.08%  '
$ws.Range('D14').Value = "'0.000009996"
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('D15').Value = "'0.6315"
$ws.Range('E15').Value = '  -6.24%  '
$ws.Range('D16').Value = '2.066.70'
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('D17').Value = "'78.59"
$ws.Range('E17').Value = '  -3.62%  '
$ws.Range('D18').Value = "'5.862"
$ws.Range('E18').Value = '  -5.66%  '
$ws.Range('D19').Value = '29.180.19'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('D20').Value = "'226.65"
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').Value = "'11.75"
$ws.Range('E22').Value = '  -4.30%  '
$ws.Range('D23').Value = "'6.991"
$ws.Range('E23').Value = '  -4.63%  '
$ws.Range('D24').Value = "'1.000"
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = "'155.80"
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('D26').Value = "'8.038"
$ws.Range('E26').Value = '  -4.99%  '
$ws.Range('E27').Value = '  -2.80%  '
$ws.Range('D28').Value = "'16.58"
$ws.Range('E28').Value = '  -4.53%  '
$ws.Range('D29').Value = "'1.493"
$ws.Range('E29').Value = '  +1.93%  '
$ws.Range('D30').Value = "'0.06327"
$ws.Range('E30').Value = '  -14.75%  '
$ws.Range('D31').Value = "'1.451"
$ws.Range('E31').Value = '  -1.79%  '
$ws.Range('D32').Value = "'3.824"
$ws.Range('E32').Value = '  -5.09%  '
$ws.Range('E33').Value = '  -5.70%  '
$ws.Range('D34').Value = "'1.123"
$ws.Range('E34').Value = '  -1.32%  '
$ws.Range('D35').Value = "'1.740"
$ws.Range('E35').Value = '  -4.33%  '
$ws.Range('D36').Value = "'0.6445"
$ws.Range('E36').Value = '  -7.08%  '
$ws.Range('D37').Value = "'2.544"
$ws.Range('E37').Value = '  -1.30%  '
$ws.Range('D38').Value = '1.214.15'
$ws.Range('E38').Value = '  -1.40%  '
$ws.Range('D39').Value = "'2.721"
$ws.Range('E39').Value = '  -3.03%  '
$ws.Range('E40').Value = '  -5.36%  '
$ws.Range('D41').Value = "'6.535"
$ws.Range('E41').Value = '  -5.02%  '
$ws.Range('D42').Value = "'0.9104"
$ws.Range('E42').Value = '  -2.14%  '
$ws.Range('D43').Value = "'1.002"
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = "'100.90"
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('D45').Value = '1.976.97'
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('E46').Value = '  -3.94%  '
$ws.Range('D47').Value = "'0.00000000118"
$ws.Range('E47').Value = '  -1.78%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = "'1.602"
$ws.Range('E48').Value = '  -5.83%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'8.516"
$ws.Range('E49').Value = '  -3.78%  '
$ws.Range('D50').Value = "'0.4568"
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('D51').Value = "'0.05517"
$ws.Range('E51').Value = '  -2.61%  '
Write-Output 'done'
